$d = $word.ActiveDocument

# Find the paragraph that ends the CREATE TABLE statement for customers
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*mobileno_UNIQUE*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find target paragraph"
}

$lines = @(
    "",
    "ALTER TABLE ``uber``.``customers`` ",
    "ADD COLUMN ``latitude`` DOUBLE NOT NULL AFTER ``approved``,",
    "ADD COLUMN ``longitude`` DOUBLE NOT NULL AFTER ``latitude``,",
    "ADD COLUMN ``rating`` DOUBLE NOT NULL AFTER ``longitude``;"
)

$curIndex = $targetIndex
foreach ($line in $lines) {
    $r = $d.Paragraphs.Item($curIndex).Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $curIndex = $curIndex + 1
    if ($line -ne "") {
        $newPara = $d.Paragraphs.Item($curIndex)
        $newPara.Range.InsertAfter($line)
    }
}
